# Initial check-in of translations changes.
#
# The "display.title" / "display.text" setting keys used throughout the
# mating_event form definition are being renamed to the more explicit
# "display.title.text" / "display.prompt.text" forms, and the active
# selection moves back onto the "survey" sheet.

$wb = $excel.ActiveWorkbook

$wsSurvey = $wb.Worksheets.Item("survey")
$wsSettings = $wb.Worksheets.Item("settings")
$wsProperties = $wb.Worksheets.Item("properties")

# Rename the translation keys referenced from the "settings" and "survey"
# sheets (order matters: it controls where the new shared strings land).
$wsSettings.Range("C1").Value = "display.title.text"
$wsSurvey.Range("F1").Value = "display.prompt.text"

# Restore the original selections on "properties" (no-op, stays E5) and
# "settings" (moves from A3 to C2).
$wsProperties.Range("E5").Select()
$wsSettings.Range("C2").Select()

# Make "survey" the active sheet/tab again, with the selection on F2.
$wsSurvey.Activate()
$wsSurvey.Range("F2").Select()
